$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three header labels whose text changed.
$ws.Range("C1").Value = "Pontosság (%)"
$ws.Range("D1").Value = "Átlagos feldolgozási idő (s)"
$ws.Range("A1").Value = "LLM modell"

# Move selection to B25, matching the saved selection state in the diff.
$ws.Range("B25").Select()
